# ============================================================
# Automatic BRVM data refresh (GitHub Actions)
# Rewrites the data rows of 'Recommandations' (A2:G49) and
# 'Top_YTD' (B2:B11) with the latest scrape.
# ============================================================

$wb = $excel.ActiveWorkbook

# --- Recommandations sheet: Titre, Jours en Hausse, Jours en Baisse,
#     Variation Totale (%), Derniere Variation (%), Recommandation, Strategie
$recoData = @(
    ,("NEI-CEDA CI", 0, 5, 4615, 930, "🟡 Observer", "➖ Neutre")
    ,("BRVM - SERVICES PUBLICS", 0, 10, 4229.5, 111.93, "🟡 Observer", "➖ Neutre")
    ,("AIR LIQUIDE CI", 0, 5, 3485, 720, "🟡 Observer", "➖ Neutre")
    ,("SUCRIVOIRE", 0, 3, 2975, 990, "🟡 Observer", "➖ Neutre")
    ,("BRVM - AUTRES SECTEURS", 0, 5, 2934.88, 585.37, "🟡 Observer", "➖ Neutre")
    ,("BRVM - DISTRIBUTION", 0, 5, 2582.81, 532.37, "🟡 Observer", "➖ Neutre")
    ,("BRVM - TRANSPORT", 0, 5, 1802.75, 359.82, "🟡 Observer", "➖ Neutre")
    ,("BRVM - AGRICULTURE", 0, 5, 1704.53, 344.59, "🟡 Observer", "➖ Neutre")
    ,("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 5, 889.23, 184.48, "🟡 Observer", "➖ Neutre")
    ,("BRVM - FINANCES", 0, 5, 726.26, 145.23, "🟡 Observer", "➖ Neutre")
    ,("BRVM - SERVICES FINANCIERS", 0, 5, 713.76, 142.73, "🟡 Observer", "➖ Neutre")
    ,("BRVM-PRESTIGE", 0, 5, 705.85, 141.01, "🟡 Observer", "➖ Neutre")
    ,("BRVM-PRINCIPAL                    (**)", 0, 3, 661.47, 221.08, "🟡 Observer", "➖ Neutre")
    ,("BRVM - INDUSTRIELS", 0, 5, 611.37, 122.49, "🟡 Observer", "➖ Neutre")
    ,("BRVM - ENERGIE", 0, 5, 563.21, 112.84, "🟡 Observer", "➖ Neutre")
    ,("BRVM - INDUSTRIE                 (**)", 0, 2, 534.34, 267.49, "🟡 Observer", "➖ Neutre")
    ,("BRVM - TELECOMMUNICATIONS", 0, 5, 477.01, 95.13, "🟡 Observer", "➖ Neutre")
    ,("BRVM - CONSOMMATION DE BASE         (**)", 0, 2, 449.09, 224.23, "🟡 Observer", "➖ Neutre")
    ,("BRVM - INDUSTRIE                  (**)", 0, 1, 266.59, 266.59, "🟡 Observer", "➖ Neutre")
    ,("BRVM - CONSOMMATION DE BASE          (**)", 0, 1, 222.65, 222.65, "🟡 Observer", "➖ Neutre")
    ,("UNILEVER CI (UNLC)", 2, 0, 14.97, 7.47, "🟡 Observer", "➖ Neutre")
    ,("CFAO MOTORS CI (CFAC)", 2, 0, 12.26, 7.28, "🟡 Observer", "➖ Neutre")
    ,("SAPH CI (SPHC)", 2, 0, 8.619999999999999, 3.92, "🟡 Observer", "➖ Neutre")
    ,("BANK OF AFRICA ML (BOAM)", 1, 0, 5.28, 5.28, "🟡 Observer", "➖ Neutre")
    ,("SERVAIR ABIDJAN CI (ABJC)", 1, 1, 4.34, 6.84, "🟡 Observer", "👀 À surveiller")
    ,("UNIWAX CI (UNXC)", 1, 1, 3.9, -3.24, "🟡 Observer", "👀 À surveiller")
    ,("SICABLE CI (CABC)", 3, 2, 3.75, 4.07, "🟡 Observer", "👀 À surveiller")
    ,("TOTALENERGIES MARKETING SN (TTLS)", 1, 0, 3.02, 3.02, "🟡 Observer", "➖ Neutre")
    ,("SOLIBRA CI (SLBC)", 1, 1, 2.22, -1.45, "🟡 Observer", "👀 À surveiller")
    ,("BERNABE CI (BNBC)", 1, 0, 2.2, 2.2, "🟡 Observer", "➖ Neutre")
    ,("ORANGE COTE D'IVOIRE (ORAC)", 1, 1, 0.68, -3.62, "🟡 Observer", "👀 À surveiller")
    ,("SICOR CI (SICC)", 1, 1, 0.5, 7.37, "🟡 Observer", "👀 À surveiller")
    ,("SOCIETE GENERALE COTE D'IVOIRE (SGBC)", 1, 1, 0.18, -4.16, "🟡 Observer", "👀 À surveiller")
    ,("BANK OF AFRICA SENEGAL (BOAS)", 1, 1, 0.13, 3.58, "🟡 Observer", "👀 À surveiller")
    ,("FILTISAC CI (FTSC)", 1, 1, 0.1, 7.43, "🟡 Observer", "👀 À surveiller")
    ,("TOTAL", 0, 2, 0, 0, "🟡 Observer", "➖ Neutre")
    ,("SETAO CI (STAC)", 2, 1, -0.02, 2.82, "🟡 Observer", "👀 À surveiller")
    ,("SAFCA CI (SAFC)", 1, 1, -0.5600000000000001, 5.02, "🟡 Observer", "👀 À surveiller")
    ,("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 1, 1, -0.5600000000000001, 3.42, "🟡 Observer", "👀 À surveiller")
    ,("NESTLE CI (NTLC)", 0, 1, -2.54, -2.54, "🟡 Observer", "➖ Neutre")
    ,("ORAGROUP TOGO (ORGT)", 0, 1, -2.77, -2.77, "🟡 Observer", "➖ Neutre")
    ,("SMB CI (SMBC)", 0, 1, -3.11, -3.11, "🟡 Observer", "➖ Neutre")
    ,("BICI CI (BICC)", 0, 1, -3.23, -3.23, "🟡 Observer", "➖ Neutre")
    ,("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -3.47, -3.47, "🟡 Observer", "➖ Neutre")
    ,("VIVO ENERGY CI (SHEC)", 1, 2, -4.12, -1.62, "🟡 Observer", "👀 À surveiller")
    ,("ECOBANK TRANS. INCORP. TG (ETIT)", 0, 1, -4.17, -4.17, "🟡 Observer", "➖ Neutre")
    ,("CIE CI (CIEC)", 0, 1, -4.58, -4.58, "🟡 Observer", "➖ Neutre")
    ,("NEI-CEDA CI (NEIC)", 0, 2, -6.42, -2.16, "🟡 Observer", "➖ Neutre")
)

# --- Top_YTD sheet: Titre, Progression YTD (%)
$ytdData = @(
    ,("BRVM - SERVICES PUBLICS", 174896530.23)
    ,("NEI-CEDA CI", 11193940)
    ,("AIR LIQUIDE CI", 3213250.4)
    ,("BRVM - AUTRES SECTEURS", 1529680.44)
    ,("BRVM - DISTRIBUTION", 890019.9300000001)
    ,("BRVM - TRANSPORT", 207071.41)
    ,("BRVM - AGRICULTURE", 166507.57)
    ,("SUCRIVOIRE", 129996.95)
    ,("BRVM - CONSOMMATION DISCRETIONNAIRE", 16440.1)
    ,("BRVM - FINANCES", 8772.799999999999)
)

$wsReco = $wb.Worksheets.Item("Recommandations")
$r = 2
foreach ($row in $recoData) {
    for ($c = 1; $c -le 7; $c++) {
        $wsReco.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}

$wsYtd = $wb.Worksheets.Item("Top_YTD")
$r = 2
foreach ($row in $ytdData) {
    for ($c = 1; $c -le 2; $c++) {
        $wsYtd.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}

Write-Host "BRVM data refreshed: Recommandations A2:G49, Top_YTD B2:B11"
